# Update the parameter listing worksheet with two new parameter rows:
#   max_waiting_time  (row 9)
#   maxvehicles       (row 10)
# and move the on-screen selection down near the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: max_waiting_time -------------------------------------------
$ws.Range("A9").Value = "max_waiting_time"
$ws.Range("B9").Value = "integer"
$ws.Range("C9").Value = "(0 - 1000) ticks "
$ws.Range("D9").Value = "The maximum amount of time passengers should wait before considering adding a new vehicle"
$ws.Range("A9:D9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 58

# --- Row 10: maxvehicles -------------------------------------------------
$ws.Range("A10").Value = "maxvehicles"
$ws.Range("B10").Value = "integer"
$ws.Range("D10").Value = "The maximum number of buses in a given route"
$ws.Range("C10").Value = "(1-10) units"
$ws.Range("A10:D10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 29

# --- Update the visible selection / scroll position ----------------------
$ws.Range("C11").Select()
